# Auto-generated script applying the crypto price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "61.069.72"
$ws.Range("E2").Value = "  +3.96%  "
# Row 3: Ethereum
$ws.Range("D3").Value = "2.720.70"
$ws.Range("E3").Value = "  +3.21%  "
# Row 4: TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.02%  "
# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "529.42"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.20%  "
# Row 6: Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.34"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.74%  "
# Row 7: USDC
$ws.Range("E7").Value = "  -0.15%  "
# Row 8: XRP
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.579"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.98%  "
# Row 9: LidoStakedEther
$ws.Range("D9").Value = "2.744.53"
$ws.Range("E9").Value = "  +3.56%  "
# Row 10: Toncoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.10"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +11.22%  "
# Row 11: Dogecoin
$ws.Range("E11").Value = "  +1.57%  "
# Row 12: Cardano
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.342"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.23%  "
# Row 13: TRON
$ws.Range("E13").Value = "  +2.79%  "
# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.202.18"
$ws.Range("E14").Value = "  +3.35%  "
# Row 15: WrappedBTC
$ws.Range("D15").Value = "61.049.64"
$ws.Range("E15").Value = "  +3.87%  "
# Row 16: Avalanche
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.56"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.56%  "
# Row 17: WrappedEther
$ws.Range("D17").Value = "2.738.89"
$ws.Range("E17").Value = "  +3.27%  "
# Row 18: ShibaInu
$ws.Range("E18").Value = "  +1.70%  "
# Row 19: BitcoinCash
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "346.51"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.45%  "
# Row 20: Polkadot
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.51"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.21%  "
# Row 21: Chainlink
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.57"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.34%  "
# Row 22: Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.41"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +4.06%  "
# Row 23: Dai
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.02%  "
# Row 24: Litecoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.44"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.79%  "
# Row 25: Kaspa
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.171"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +4.71%  "
# Row 26: Polygon
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.420"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.49%  "
# Row 27: Binance-PegBSC-USD
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.994"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.29%  "
# Row 28: PEPE
$ws.Range("D28").Value = "0.0₃0828"
$ws.Range("E28").Value = "  +3.08%  "
# Row 29: InternetComputer(DFINITY)
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.34"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.50%  "
# Row 30: Aptos
$ws.Range("E30").Value = "  +8.18%  "
# Row 31: USDe
$ws.Range("E31").Value = "  -0.05%  "
# Row 32: PancakeSwap
$ws.Range("E32").Value = "  +2.28%  "
# Row 33: EthereumClassic
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.11"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.23%  "
# Row 34: Monero
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.42"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.01%  "
# Row 35: NEARProtocol
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.26"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +7.03%  "
# Row 36: ImmutableX
$ws.Range("E36").Value = "  +7.72%  "
# Row 37: SuiNetwork
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.924"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -7.00%  "
# Row 38: Fetch.AI
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.897"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +7.03%  "
# Row 39: Stacks
$ws.Range("E39").Value = "  +8.98%  "
# Row 40: OKB
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.73"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.21%  "
# Row 41: Filecoin
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.68"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.39%  "
# Row 42: Mantle
$ws.Range("E42").Value = "  +4.86%  "
# Row 43: EnergySwap
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.28"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.41%  "
# Row 44: Bittensor
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "280.89"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.63%  "
# Row 45: FirstDigitalUSD
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.995"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.36%  "
# Row 46: Stellar
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0985"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.19%  "
# Row 47: Maker/RenderToken (swap)
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.99"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +7.65%  "
# Row 48: RenderToken/Maker (swap)
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "2.102.00"
$ws.Range("E48").Value = "  +5.52%  "
# Row 49: Hedera
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0542"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.66%  "
# Row 50: InjectiveProtocol/WhiteBITCoin (swap)
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.53"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.18%  "
# Row 51: WhiteBITCoin/InjectiveProtocol (swap)
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.45"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +5.36%  "
